$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values: C2 and E2 are cleared (removed), B2 and D2 updated
$ws.Range("B2").Value = 17.375431314037488
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 26.021814370035074
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 14.471734435433772
$ws.Range("C3").Value = -10.616310651571711
$ws.Range("D3").Value = 25.673148807290115
$ws.Range("E3").Value = -5.1521007428528565

# Update the selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
